$p = $ppt.ActivePresentation

# --- Slide 1: "Today's date" placeholder -------------------------------
# Originally built from five runs: "May" + " " + "5" + ", " + "2016".
# The first two runs ("May" and " ") get merged into a single run "May "
# (this leaves the visible text "May 5, 2016" unchanged, only the run
# split changes), so we rewrite just that leading span of characters.
$slide1 = $p.Slides.Item(1)
$dateShape = $slide1.Shapes.Item(2)
$dateRange = $dateShape.TextFrame.TextRange
$dateRange.Characters(1, 4).Text = "May "

# --- Slide 11: "Alert" -> "Prompt" --------------------------------------
$slide11 = $p.Slides.Item(11)
$contentShape = $slide11.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange
$contentRange.Characters(39, 6).Text = "Prompt "
